# Update stats for 2025-11 (row 24)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C24").Value = 999
$ws.Range("D24").Value = 5950730
$ws.Range("E24").Value = 936.8277707808564
$ws.Range("G24").Value = 3.523316062176174
$ws.Range("H24").Value = 26.05951143208916
